$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "52-34="
$t.Cell(1,2).Range.Text = "42+34="
$t.Cell(1,3).Range.Text = "49+47="
$t.Cell(1,4).Range.Text = "51-20="
$t.Cell(1,5).Range.Text = "51-16="
$t.Cell(2,1).Range.Text = "75+12="
$t.Cell(2,2).Range.Text = "77-59="
$t.Cell(2,3).Range.Text = "3+90="
$t.Cell(2,4).Range.Text = "70-8="
$t.Cell(2,5).Range.Text = "59+22="
$t.Cell(3,1).Range.Text = "40-30="
$t.Cell(3,2).Range.Text = "47+27="
$t.Cell(3,3).Range.Text = "51+43="
$t.Cell(3,4).Range.Text = "35+32="
$t.Cell(3,5).Range.Text = "67+29="
$t.Cell(4,1).Range.Text = "26+3="
$t.Cell(4,2).Range.Text = "58-36="
$t.Cell(4,3).Range.Text = "66-50="
$t.Cell(4,4).Range.Text = "74-30="
$t.Cell(4,5).Range.Text = "70-15="
$t.Cell(5,1).Range.Text = "35-1="
$t.Cell(5,2).Range.Text = "56+37="
$t.Cell(5,3).Range.Text = "97-72="
$t.Cell(5,4).Range.Text = "99-23="
$t.Cell(5,5).Range.Text = "84-31="
$t.Cell(6,1).Range.Text = "67-12="
$t.Cell(6,2).Range.Text = "56-51="
$t.Cell(6,3).Range.Text = "35+19="
$t.Cell(6,4).Range.Text = "76-43="
$t.Cell(6,5).Range.Text = "10-9="
$t.Cell(7,1).Range.Text = "62+7="
$t.Cell(7,2).Range.Text = "44-34="
$t.Cell(7,3).Range.Text = "33+18="
$t.Cell(7,4).Range.Text = "24+10="
$t.Cell(7,5).Range.Text = "79-1="
$t.Cell(8,1).Range.Text = "52+27="
$t.Cell(8,2).Range.Text = "67-0="
$t.Cell(8,3).Range.Text = "19+4="
$t.Cell(8,4).Range.Text = "29+15="
$t.Cell(8,5).Range.Text = "70-17="
$t.Cell(9,1).Range.Text = "2-0="
$t.Cell(9,2).Range.Text = "42+40="
$t.Cell(9,3).Range.Text = "66-63="
$t.Cell(9,4).Range.Text = "63-36="
$t.Cell(9,5).Range.Text = "72-50="
$t.Cell(10,1).Range.Text = "96-2="
$t.Cell(10,2).Range.Text = "26-22="
$t.Cell(10,3).Range.Text = "12+3="
$t.Cell(10,4).Range.Text = "60+38="
$t.Cell(10,5).Range.Text = "1+36="
$t.Cell(11,1).Range.Text = "97-64="
$t.Cell(11,2).Range.Text = "42+30="
$t.Cell(11,3).Range.Text = "54-51="
$t.Cell(11,4).Range.Text = "62-12="
$t.Cell(11,5).Range.Text = "22+21="
$t.Cell(12,1).Range.Text = "38+12="
$t.Cell(12,2).Range.Text = "49+31="
$t.Cell(12,3).Range.Text = "29+39="
$t.Cell(12,4).Range.Text = "89-58="
$t.Cell(12,5).Range.Text = "35+45="
$t.Cell(13,1).Range.Text = "38-31="
$t.Cell(13,2).Range.Text = "73-29="
$t.Cell(13,3).Range.Text = "39+29="
$t.Cell(13,4).Range.Text = "4+62="
$t.Cell(13,5).Range.Text = "78-45="
$t.Cell(14,1).Range.Text = "46-46="
$t.Cell(14,2).Range.Text = "42-34="
$t.Cell(14,3).Range.Text = "23+49="
$t.Cell(14,4).Range.Text = "73-66="
$t.Cell(14,5).Range.Text = "10+74="
$t.Cell(15,1).Range.Text = "7+69="
$t.Cell(15,2).Range.Text = "71-70="
$t.Cell(15,3).Range.Text = "87-56="
$t.Cell(15,4).Range.Text = "43-11="
$t.Cell(15,5).Range.Text = "32-5="
$t.Cell(16,1).Range.Text = "80-75="
$t.Cell(16,2).Range.Text = "83-16="
$t.Cell(16,3).Range.Text = "48+44="
$t.Cell(16,4).Range.Text = "51+38="
$t.Cell(16,5).Range.Text = "59-7="
$t.Cell(17,1).Range.Text = "62+21="
$t.Cell(17,2).Range.Text = "14+27="
$t.Cell(17,3).Range.Text = "12+65="
$t.Cell(17,4).Range.Text = "84-60="
$t.Cell(17,5).Range.Text = "50-18="
$t.Cell(18,1).Range.Text = "50+19="
$t.Cell(18,2).Range.Text = "37+16="
$t.Cell(18,3).Range.Text = "67-60="
$t.Cell(18,4).Range.Text = "10+81="
$t.Cell(18,5).Range.Text = "62+31="
$t.Cell(19,1).Range.Text = "17+82="
$t.Cell(19,2).Range.Text = "42+27="
$t.Cell(19,3).Range.Text = "30+31="
$t.Cell(19,4).Range.Text = "93-51="
$t.Cell(19,5).Range.Text = "22+4="
$t.Cell(20,1).Range.Text = "30-29="
$t.Cell(20,2).Range.Text = "8+59="
$t.Cell(20,3).Range.Text = "81-24="
$t.Cell(20,4).Range.Text = "60-34="
$t.Cell(20,5).Range.Text = "19-17="
